$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quellen")
$ws.Rows("12:12").Delete()
$ws.Activate()
$ws.Range("B22").Select()
